# "working containers and sites parsers"
# The "Crew" column is removed from each of the four worksheets (Init, Picking,
# Shocking, HU Transfer). Deleting the entire column also removes the header
# comment that was attached to it ("Eg. AB, CD") and shifts every column to its
# right (e.g. "Comments") one position to the left. Finally, the "Picking"
# worksheet becomes the active/selected sheet in the workbook.

$wb = $excel.ActiveWorkbook

# --- Init sheet: Crew column is I ---
$ws1 = $wb.Worksheets.Item("Init")
$ws1.Range("I3").Comment.Delete()
$ws1.Columns("I").Delete()
$ws1.Columns("I").Select()

# --- Picking sheet: Crew column is J ---
$ws2 = $wb.Worksheets.Item("Picking")
$ws2.Range("J3").Comment.Delete()
$ws2.Columns("J").Delete()
$ws2.Columns("J").Select()

# --- Shocking sheet: Crew column is J ---
$ws3 = $wb.Worksheets.Item("Shocking")
$ws3.Range("J3").Comment.Delete()
$ws3.Columns("J").Delete()
$ws3.Columns("J").Select()

# --- HU Transfer sheet: Crew column is Q ---
$ws4 = $wb.Worksheets.Item("HU Transfer")
$ws4.Range("Q3").Comment.Delete()
$ws4.Columns("Q").Delete()
$ws4.Columns("Q").Select()

# Picking becomes the active tab (activeTab index 1, 0-based)
$ws2.Activate()
